$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells are written as text (avoid Excel auto-converting
# strings like "1.00" or "33.40" into numbers and losing the original formatting).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.005.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.945.58"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.35%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "375.56"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.25%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.28%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.64%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.37"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.33%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.83%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.403.40"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.58"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "11.36"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +52.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.948.42"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.00"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.89%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "50.973.18"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.08"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.48"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "265.65"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.88"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.14"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.15"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.52"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.12"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.46%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.28%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.64%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.02"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.95"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "33.40"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.14%  "

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Toncoin"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.02"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.17%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.21%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.17"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.14%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.69%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.81"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.48"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.59"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.27"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.24%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.41"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.00%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.81%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.62%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.75%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.992.62"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0327"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.97%  "
